$wb = $excel.ActiveWorkbook

# Insert a new worksheet "emplVerification" right after "validLoginData"
# (i.e. before "invalidLoginData"), matching the diff which adds it as the
# new second sheet with sheetId=3.
$validSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $validSheet)
$newSheet.Name = "emplVerification"

# Populate the new sheet's data
$newSheet.Range("A1").Value = "emp_id"
$newSheet.Range("B1").Value = "emp_name"
$newSheet.Range("A2").Value = 2
$newSheet.Range("B2").Value = "yura"

# Make this the active/selected sheet with B3 as the active cell/selection
[void]$newSheet.Range("B3").Select()
